$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Emails that should be moved to the front of the "Recorded By" list
# (admin@admin.com is intentionally excluded - it stays put)
$movable = @("backup@backdoor.com", "dnasr281@gmail.com")

$used = $ws.UsedRange
$firstRow = $used.Row
$lastRow = $firstRow + $used.Rows.Count - 1

for ($r = $firstRow; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $text = $cell.Text

    if ($text -eq $null -or $text -eq "") {
        continue
    }
    if ($text.IndexOf(",") -lt 0) {
        continue
    }

    $rawParts = $text.Split(",")
    $parts = @()
    foreach ($p in $rawParts) {
        $parts += $p.Trim()
    }

    $front = @()
    $back = @()
    foreach ($p in $parts) {
        if ($movable -contains $p) {
            $front += $p
        } else {
            $back += $p
        }
    }

    if ($front.Length -eq 0) {
        continue
    }

    $newParts = $front + $back
    $newText = $newParts -join ", "

    if ($newText -ne $text) {
        $cell.Value = $newText
    }
}
